$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34; existing rows 34-50 shift down to 35-51.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly price record.
$ws.Cells.Item(34, 1).Value = 4
$ws.Cells.Item(34, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(34, 3).Value = "Los Lagos"
$ws.Cells.Item(34, 4).Value = 44510
$ws.Cells.Item(34, 5).Value = 10
$ws.Cells.Item(34, 6).Value = 100112031
$ws.Cells.Item(34, 7).Value = "Poroto verde"
$ws.Cells.Item(34, 8).Value = "Magnum"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 40
$ws.Cells.Item(34, 11).Value = 40000
$ws.Cells.Item(34, 12).Value = 40000
$ws.Cells.Item(34, 13).Value = 40000
$ws.Cells.Item(34, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(34, 15).Value = "Perú"
$ws.Cells.Item(34, 16).Value = 1600
$ws.Cells.Item(34, 17).Value = 25
$ws.Cells.Item(34, 18).Value = "Hortaliza"
